$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 760.2222
$ws.Range("I2").Value = 193.4
$ws.Range("K2").Value = 193.4
$ws.Range("M2").Value = -80.40000000000001
$ws.Range("H28").Value = 2630.15
$ws.Range("I28").Value = 2056.4
$ws.Range("J28").Value = 3203.9
$ws.Range("K28").Value = 2056.4
$ws.Range("L28").Value = 3203.9
$ws.Range("M28").Value = -1571.4
$ws.Range("N28").Value = -4173.9
$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877
$ws.Range("H88").Value = 1504.2
$ws.Range("I88").Value = 1110.6
$ws.Range("J88").Value = 1897.8
$ws.Range("K88").Value = 1110.6
$ws.Range("L88").Value = 1897.8
$ws.Range("M88").Value = -704.5999999999999
$ws.Range("N88").Value = -2709.8
$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384
$ws.Range("H91").Value = 1504.2
$ws.Range("I91").Value = 1110.6
$ws.Range("J91").Value = 1897.8
$ws.Range("K91").Value = 1110.6
$ws.Range("L91").Value = 1897.8
$ws.Range("M91").Value = 293.4000000000001
$ws.Range("N91").Value = -4705.8
$ws.Range("H103").Value = 639
$ws.Range("I103").Value = 621
$ws.Range("J103").Value = 663
$ws.Range("K103").Value = 1863
$ws.Range("L103").Value = 1989
$ws.Range("M103").Value = -1277
$ws.Range("N103").Value = -3161
$ws.Range("H107").Value = 974.7143
$ws.Range("I107").Value = 591.4545000000001
$ws.Range("J107").Value = 2380
$ws.Range("K107").Value = 591.4545000000001
$ws.Range("L107").Value = 2380
$ws.Range("M107").Value = 1328.5455
$ws.Range("N107").Value = -6220
$ws.Range("H116").Value = 3802
$ws.Range("J116").Value = 3498.5
$ws.Range("L116").Value = 3498.5
$ws.Range("N116").Value = -10382.5
$ws.Range("H137").Value = 2944.077
$ws.Range("I137").Value = 2346.6897
$ws.Range("J137").Value = 3425.3057
$ws.Range("K137").Value = 7040.0691
$ws.Range("L137").Value = 10275.9171
$ws.Range("M137").Value = -4490.0691
$ws.Range("N137").Value = -15375.9171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4862.6807
$ws.Range("I32").Value = 2589.7144
$ws.Range("K32").Value = 2589.7144
$ws.Range("M32").Value = -2302.7144
$ws.Range("H92").Value = 58775
$ws.Range("J92").Value = 58775
$ws.Range("L92").Value = 58775
$ws.Range("N92").Value = -63767
$ws.Range("H132").Value = 2059.75
$ws.Range("I132").Value = 2059.75
$ws.Range("K132").Value = 6179.25
$ws.Range("M132").Value = -3649.25
$ws.Range("H135").Value = 52785.668
$ws.Range("J135").Value = 52785.668
$ws.Range("L135").Value = 52785.668
$ws.Range("N135").Value = -62925.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 690
$ws.Range("I20").Value = 690
$ws.Range("K20").Value = 690
$ws.Range("M20").Value = -443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 495.08334
$ws.Range("I22").Value = 503.72726
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 503.72726
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -153.72726
$ws.Range("N22").Value = -1100
$ws.Range("H45").Value = 70
$ws.Range("I45").Value = 70
$ws.Range("K45").Value = 70
$ws.Range("M45").Value = 523
$ws.Range("H86").Value = 2377.8333
$ws.Range("I86").Value = 2040.6666
$ws.Range("K86").Value = 2040.6666
$ws.Range("M86").Value = -917.6666
$ws.Range("H89").Value = 2377.8333
$ws.Range("I89").Value = 2040.6666
$ws.Range("K89").Value = 10203.333
$ws.Range("M89").Value = -4587.333000000001
$ws.Range("H99").Value = 2935.4
$ws.Range("I99").Value = 2247.6667
$ws.Range("K99").Value = 2247.6667
$ws.Range("M99").Value = -749.6667000000002
$ws.Range("H126").Value = 2935.4
$ws.Range("I126").Value = 2247.6667
$ws.Range("K126").Value = 6743.000100000001
$ws.Range("M126").Value = -4273.000100000001
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 283
$ws.Range("I5").Value = 279.6
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 838.8000000000001
$ws.Range("L5").Value = 900
$ws.Range("M5").Value = -726.8000000000001
$ws.Range("N5").Value = -1124
$ws.Range("H14").Value = 262.54544
$ws.Range("I14").Value = 262.54544
$ws.Range("K14").Value = 787.63632
$ws.Range("M14").Value = -614.63632
$ws.Range("H68").Value = 2779.2
$ws.Range("I68").Value = 974.25
$ws.Range("K68").Value = 2922.75
$ws.Range("M68").Value = -2111.75
$ws.Range("H71").Value = 2779.2
$ws.Range("I71").Value = 974.25
$ws.Range("K71").Value = 8768.25
$ws.Range("M71").Value = -4712.25
$ws.Range("H92").Value = 173.83333
$ws.Range("I92").Value = 117.5
$ws.Range("K92").Value = 352.5
$ws.Range("M92").Value = 895.5
$ws.Range("H107").Value = 1069.6471
$ws.Range("I107").Value = 725.75
$ws.Range("J107").Value = 1175.4615
$ws.Range("K107").Value = 2177.25
$ws.Range("L107").Value = 3526.3845
$ws.Range("M107").Value = -257.25
$ws.Range("N107").Value = -7366.3845
$ws.Range("H131").Value = 28151448
$ws.Range("I131").Value = 1111111
$ws.Range("K131").Value = 3333333
$ws.Range("M131").Value = -3328293
$ws.Range("H135").Value = 283
$ws.Range("I135").Value = 279.6
$ws.Range("J135").Value = 300
$ws.Range("K135").Value = 2516.4
$ws.Range("L135").Value = 2700
$ws.Range("M135").Value = 18.59999999999991
$ws.Range("N135").Value = -7770

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1657.4
$ws.Range("I113").Value = 1444
$ws.Range("K113").Value = 1444
$ws.Range("M113").Value = 726
$ws.Range("H132").Value = 1851.8334
$ws.Range("I132").Value = 1851.8334
$ws.Range("K132").Value = 5555.5002
$ws.Range("M132").Value = -3025.5002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3581.1667
$ws.Range("I7").Value = 2574.25
$ws.Range("J7").Value = 4084.625
$ws.Range("K7").Value = 2574.25
$ws.Range("L7").Value = 4084.625
$ws.Range("M7").Value = -2462.25
$ws.Range("N7").Value = -4308.625
$ws.Range("H74").Value = 28333
$ws.Range("I74").Value = 28333
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 28333
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -27335
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 28333
$ws.Range("I77").Value = 28333
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 84999
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -80007
$ws.Range("N77").ClearContents()
$ws.Range("H126").Value = 3581.1667
$ws.Range("I126").Value = 2574.25
$ws.Range("J126").Value = 4084.625
$ws.Range("K126").Value = 7722.75
$ws.Range("L126").Value = 12253.875
$ws.Range("M126").Value = -5252.75
$ws.Range("N126").Value = -17193.875
$ws.Range("H136").Value = 2831.077
$ws.Range("I136").Value = 2749.0908
$ws.Range("J136").Value = 3282
$ws.Range("K136").Value = 8247.2724
$ws.Range("L136").Value = 9846
$ws.Range("M136").Value = -5697.2724
$ws.Range("N136").Value = -14946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3242.5
$ws.Range("I126").Value = 2614.75
$ws.Range("K126").Value = 7844.25
$ws.Range("M126").Value = -5374.25
$ws.Range("H136").Value = 1400.0454
$ws.Range("I136").Value = 1175.4375
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 3526.3125
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -976.3125
$ws.Range("N136").Value = -11097
$ws.Range("H141").Value = 67500
$ws.Range("J141").Value = 67500
$ws.Range("L141").Value = 67500
$ws.Range("N141").Value = -77860
